$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 44: "Ni3(HAB)2" / 2D layered / BET 133 / Porous Yes / Year 2021 / Value 0.52 / DOI
# Write the DOI (column G) first so the shared-string table gains the new DOI
# entry before the new compound name, matching the authored order of the
# shared strings table (DOI string lands at index 58, name at index 59).
$ws.Range("G44").Value = "10.1021/acsaem.0c02758 "
$ws.Range("A44").Value = "Ni3(HAB)2"
$ws.Range("B44").Value = "2D layered"
$ws.Range("C44").Value = 133
$ws.Range("D44").Value = "Yes"
$ws.Range("D44").NumberFormat = "0.00E+00"
$ws.Range("E44").Value = 2021
$ws.Range("F44").Value = 0.52
$ws.Range("F44").NumberFormat = "0.00E+00"

# Move the active selection onto the newly added row, same as the source
# workbook (selection moved from G47 to G44 once the table only runs to
# row 44).
$null = $ws.Range("G44").Select()

# Best-effort: scroll the view so row 12 / column B is the top-left visible
# cell (mirrors topLeftCell="B12" on the saved sheetView).
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 2
